$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30-58 down to 31-59
$ws.Rows("30").Insert()

# Populate the new row 30 with the "Concerto inaugurale della nuova Filarmonica
# del Teatro comunale" event data
$ws.Range("A30").Value = "Concerti,Spettacoli,Musica"
$ws.Range("B30").Value = "Modena"
$ws.Range("C30").Value = "corso Canalgrande, 85"
$ws.Range("D30").Value = "2022-05-30T11:40:24+00:00"
$ws.Range("F30").Value = "2022-05-30T11:41:23+00:00"
$ws.Range("H30").Value = "2022-06-06T11:00:00+00:00"
$ws.Range("I30").Value = "2022-06-06T12:00:00+00:00"
$ws.Range("J30").Value = "https://www.comune.modena.it/api/novita/eventi/2022/concerto-inaugurale-della-nuova-filarmonica-del-teatro-comunale/@@images/bd0aa52a-489c-4715-b957-8a4551ed5f95.jpeg"
$ws.Range("K30").Value = "Dmitry Masleev"
$ws.Range("L30").Value = "2022-06-01T12:46:07+00:00"
$ws.Range("M30").Value = "Teatro comunale Pavarotti-Freni"
$ws.Range("N30").Value = " ore 21.00"
$ws.Range("P30").Value = " Ingressio con biglietto gratuito. I biglietti sono disponibili presso la biglietteria del Teatro o telefonando allo 059  2033010"
$ws.Range("S30").Value = "Concerto inaugurale della nuova Filarmonica del Teatro comunale"
$ws.Range("V30").Value = $false
$ws.Range("W30").Value = 41123
$ws.Range("X30").Value = "https://www.comune.modena.it/novita/eventi/2022/concerto-inaugurale-della-nuova-filarmonica-del-teatro-comunale"
$ws.Range("Y30").Value = "44,64582"
$ws.Range("Z30").Value = "10,92572"
$ws.Range("AA30").Value = "POINT (10.92572 44.64582)"
